$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) counts for column G, rows 2-35 (replacing previous
# placeholder "Strike#" values with regenerated K values).
$newK = @{
    2  = 2
    3  = 9
    4  = 6
    5  = 7
    6  = 6
    7  = 8
    8  = 4
    9  = 7
    10 = 4
    11 = 5
    12 = 9
    13 = 6
    14 = 7
    15 = 4
    16 = 4
    17 = 4
    18 = 6
    19 = 3
    20 = 13
    21 = 6
    22 = 9
    23 = 2
    24 = 5
    25 = 7
    26 = 7
    27 = 6
    28 = 2
    29 = 4
    30 = 5
    31 = 1
    32 = 1
    33 = 4
    34 = 2
    35 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
